$wb = $excel.ActiveWorkbook

# --- Settings sheet ---
$settings = $wb.Worksheets.Item("Settings")

# Row 2: rename the UploadBucket asset reference (drop the "[Dev] " prefix)
$settings.Range("B2").Value = "RPA_Moon_UploadBucket"

# Row 8 (TemplateMolpayFee / Data\Template\Molpay Fee.xlsx) is no longer used - clear it out
$settings.Range("A8:B8").ClearContents()

# --- Assets sheet ---
$assets = $wb.Worksheets.Item("Assets")

# Rename the "[Dev] RPA_Moon_*" asset names to their production "RPA_Moon_*" equivalents
$assets.Range("B2").Value = "RPA_Moon_SheetIdConfig"
$assets.Range("B3").Value = "RPA_Moon_PathMasterFolder"
$assets.Range("B4").Value = "RPA_Moon_PathMailTemplate"
$assets.Range("B5").Value = "RPA_Moon_PathSaKey"

# Switch the MOLPAY captcha / merchant asset names from the old RPA044/Moon naming to RPA021_MOLPAY
$assets.Range("B6").Value = "RPA021_MOLPAY_Captcha_SiteKey"
$assets.Range("B7").Value = "RPA021_MOLPAY_Captcha_RuleId"
$assets.Range("B8").Value = "RPA021_MOLPAY_MerchantId"

$assets.Range("B9").Value = "RPA_Moon_PathDownloadChrome"
$assets.Range("B10").Value = "RPA_Moon_DialogDownloadChrome"

# --- Selections / active sheet ---
$settings.Activate()
[void]$settings.Range("B3").Select()

$assets.Activate()
[void]$assets.Range("A11:B11").Select()
